$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H2 and H4 source values
$ws.Range("H2").Value = 6
$ws.Range("H4").Value = 23

# Row 10 headers (mirrors columns Lectures/Hours/Weight -> G1/H1... actually Lectures,Hours,Weight labels)
$ws.Range("C10").Value = "Lectures"
$ws.Range("D10").Value = "Hours"
$ws.Range("E10").Value = "Weight"

# Row 11 - Christina (entered individually, not part of the fill-down group)
$ws.Range("B11").Formula = "=B2"
$ws.Range("C11").Formula = "=H2"
$ws.Range("D11").Formula = "=J2"
$ws.Range("E11").Formula = "=K2"

# Rows 12:14 - Martin/Jon/Total, filled down from row 12 so Excel records them
# as a shared formula group (matches fill-down / autofill behaviour)
$ws.Range("B12:B14").Formula = "=B3"
$ws.Range("C12:C14").Formula = "=H3"

$ws.Range("D12").Formula = "=J3"
$ws.Range("E12").Formula = "=K3"
$ws.Range("D13").Formula = "=J4"
$ws.Range("E13").Formula = "=K4"
$ws.Range("D14").Formula = "=J5"
$ws.Range("E14").Formula = "=K5"

# Apply number formats matching J/K columns (percentage / decimal)
$ws.Range("D11:D14").NumberFormat = "0.0"
$ws.Range("E11:E14").NumberFormat = "0.0%"

# Update selection to match diff
$ws.Range("B10:E14").Select()

$wb.Save()
